$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A33").Value = "Angola"
Write-Host "done"
